# Digits Kinematics Data - append "(m)" unit suffix to the measurement
# column headers (A1:M1) and widen the affected columns (B:M) to fit the
# longer header text, mirroring an Excel AutoFit after the header edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append " (m)" to every header cell in row 1 (A1:M1).
for ($c = 1; $c -le 13; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value() = $cell.Value() + " (m)"
}

# 2) Re-fit column widths B:M for the now-longer header labels (column A's
#    header text, "DP_Length", did not change length so its width is left
#    alone).
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 13
$ws.Columns.Item(4).ColumnWidth = 13.5
$ws.Columns.Item(5).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 19
$ws.Columns.Item(7).ColumnWidth = 19
$ws.Columns.Item(8).ColumnWidth = 20
$ws.Columns.Item(9).ColumnWidth = 18
$ws.Columns.Item(10).ColumnWidth = 19
$ws.Columns.Item(11).ColumnWidth = 19
$ws.Columns.Item(12).ColumnWidth = 20
$ws.Columns.Item(13).ColumnWidth = 27.5
